$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '45.848.28'
$ws.Range('E2').Value = '  -1.97%  '
$ws.Range('D3').Value = '2.480.17'
$ws.Range('E3').Value = '  +9.82%  '
$ws.Range('E4').Value = '  +0.06%  '
$ws.Range('D5').Value = '''295.08'
$ws.Range('E5').Value = '  -1.03%  '
$ws.Range('D6').Value = '''95.04'
$ws.Range('E6').Value = '  -2.86%  '
$ws.Range('D7').Value = '''0.574'
$ws.Range('E7').Value = '  +0.54%  '
$ws.Range('E8').Value = '  -0.01%  '
$ws.Range('D9').Value = '''0.520'
$ws.Range('E9').Value = '  +3.69%  '
$ws.Range('D10').Value = '''35.25'
$ws.Range('E10').Value = '  +2.36%  '
$ws.Range('D11').Value = '''0.0786'
$ws.Range('E11').Value = '  -0.25%  '
$ws.Range('D12').Value = '''7.32'
$ws.Range('E12').Value = '  +5.00%  '
$ws.Range('E13').Value = '  +1.82%  '
$ws.Range('D14').Value = '2.839.09'
$ws.Range('E14').Value = '  +9.13%  '
$ws.Range('D15').Value = '2.454.56'
$ws.Range('E15').Value = '  +8.79%  '
$ws.Range('D16').Value = '''0.854'
$ws.Range('E16').Value = '  +8.06%  '
$ws.Range('D17').Value = '''14.24'
$ws.Range('E17').Value = '  +5.58%  '
$ws.Range('D18').Value = '45.890.47'
$ws.Range('E18').Value = '  -1.77%  '
$ws.Range('D19').Value = '''12.80'
$ws.Range('E19').Value = '  +4.35%  '
$ws.Range('D20').Value = '0.0₃0948'
$ws.Range('E20').Value = '  -2.29%  '
$ws.Range('D21').Value = '''6.33'
$ws.Range('E21').Value = '  +10.11%  '
$ws.Range('D22').Value = '''67.57'
$ws.Range('E22').Value = '  +3.11%  '
$ws.Range('D23').Value = '''246.71'
$ws.Range('E23').Value = '  +1.37%  '
$ws.Range('D24').Value = '''2.79'
$ws.Range('E24').Value = '  +0.85%  '
$ws.Range('B25').Value = 'ImmutableX'
$ws.Range('C25').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D25').Value = '''1.96'
$ws.Range('E25').Value = '  +6.44%  '
$ws.Range('B26').Value = 'Dai'
$ws.Range('C26').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D26').Value = '''1.00'
$ws.Range('E26').Value = '  +0.11%  '
$ws.Range('D27').Value = '''39.22'
$ws.Range('E27').Value = '  -3.98%  '
$ws.Range('D28').Value = '''2.22'
$ws.Range('E28').Value = '  +0.42%  '
$ws.Range('D29').Value = '''9.85'
$ws.Range('E29').Value = '  +4.34%  '
$ws.Range('D30').Value = '''21.91'
$ws.Range('E30').Value = '  +9.93%  '
$ws.Range('D31').Value = '''3.86'
$ws.Range('E31').Value = '  +16.70%  '
$ws.Range('D32').Value = '''2.77'
$ws.Range('E32').Value = '  -1.90%  '
$ws.Range('D33').Value = '''5.56'
$ws.Range('E33').Value = '  +5.60%  '
$ws.Range('D34').Value = '''147.58'
$ws.Range('E34').Value = '  +3.12%  '
$ws.Range('D35').Value = '''2.09'
$ws.Range('E35').Value = '  +27.02%  '
$ws.Range('D36').Value = '''0.0774'
$ws.Range('E36').Value = '  +1.63%  '
$ws.Range('D37').Value = '''0.115'
$ws.Range('E37').Value = '  +4.06%  '
$ws.Range('D38').Value = '''0.116'
$ws.Range('E38').Value = '  +0.95%  '
$ws.Range('D39').Value = '''15.34'
$ws.Range('E39').Value = '  +1.01%  '
$ws.Range('D40').Value = '''3.99'
$ws.Range('E40').Value = '  +5.76%  '
$ws.Range('D41').Value = '''0.0300'
$ws.Range('E41').Value = '  +2.46%  '
$ws.Range('D42').Value = '2.020.45'
$ws.Range('E42').Value = '  +13.52%  '
$ws.Range('D43').Value = '''3.26'
$ws.Range('E43').Value = '  +6.48%  '
$ws.Range('D44').Value = '''1.00'
$ws.Range('E44').Value = '  +0.13%  '
$ws.Range('D45').Value = '''92.16'
$ws.Range('E45').Value = '  -0.43%  '
$ws.Range('D46').Value = '''1.78'
$ws.Range('E46').Value = '  -3.89%  '
$ws.Range('D47').Value = '''16.13'
$ws.Range('B48').Value = 'Aave'
$ws.Range('C48').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D48').Value = '''103.49'
$ws.Range('E48').Value = '  +10.90%  '
$ws.Range('B49').Value = 'FraxShare'
$ws.Range('C49').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D49').Value = '''8.63'
$ws.Range('E49').Value = '  +11.21%  '
$ws.Range('D50').Value = '2.691.05'
$ws.Range('E50').Value = '  +8.59%  '
$ws.Range('D51').Value = '''0.186'
$ws.Range('E51').Value = '  +2.33%  '
